$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell value updates per diff (cryptos list refresh).
# Column D holds price strings that must remain text (e.g. "7.00", "0.999");
# Excels Value setter auto-coerces plain-looking numbers, so we briefly force
# a text number format, assign, then restore the cell to the Normal style so
# no stray style index is left on the cell.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "61.150.71"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.81%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.928.66"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.97%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "592.18"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.41%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "145.39"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.07%  "
$ws.Range("E8").Value = "  +0.69%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "7.00"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +4.89%  "
$ws.Range("E10").Value = "  +0.12%  "
$ws.Range("E11").Value = "  -0.47%  "
$ws.Range("E12").Value = "  +0.54%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "33.78"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.60%  "
$ws.Range("E14").Value = "  -0.35%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.413.54"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.09%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "61.111.20"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "6.73"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.85%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.932.10"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.10%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "436.26"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.12%  "
$ws.Range("E20").Value = "  -0.38%  "
$ws.Range("E21").Value = "  -0.66%  "
$ws.Range("E22").Value = "  +0.85%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "81.47"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.81%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.07"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.26%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.22"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.86%  "
$ws.Range("E26").Value = "  +0.34%  "
$ws.Range("E27").Value = "  +0.02%  "
$ws.Range("E28").Value = "  +2.16%  "
$ws.Range("E29").Value = "  +0.07%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.00"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.20%  "
$ws.Range("E31").Value = "  +4.69%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "26.69"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.21%  "
$ws.Range("E33").Value = "  +0.02%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0₃0869"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.69%  "
$ws.Range("E35").Value = "  -0.13%  "
$ws.Range("E36").Value = "  +0.96%  "
$ws.Range("E37").Value = "  -0.04%  "
$ws.Range("B38").Value = "Kaspa"
$ws.Range("C38").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.124"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.26%  "
$ws.Range("B39").Value = "Stacks"
$ws.Range("C39").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.00"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.95%  "
$ws.Range("E40").Value = "  +0.37%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "42.16"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +4.00%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.286"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.39%  "
$ws.Range("E43").Value = "  +1.48%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0347"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.25%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.689.61"
$ws.Range("D45").Style = "Normal"
$ws.Range("E46").Value = "  +1.15%  "
$ws.Range("E47").Value = "  -0.03%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "24.03"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.45%  "
$ws.Range("E49").Value = "  -0.33%  "
$ws.Range("E50").Value = "  -1.43%  "
$ws.Range("E51").Value = "  +0.00%  "
